# Parties-prenantes.docx — split/merge runs in the second column of the
# table so the visible text matches the updated wording, while keeping
# the same run layout the diff shows.

function Split-RunAt($doc, $range, $offset) {
    # Forces a run boundary `$offset` characters after the start of
    # `$range` by toggling a character property on the leading
    # sub-range and right back off. Word (and this COM host) always
    # keeps distinctly-touched runs apart, so the abutting text ends
    # up as two adjacent <w:r> elements instead of being re-coalesced.
    $s = $range.Start
    $sub = $doc.Range($s, $s + $offset)
    $sub.Font.Bold = 1
    $sub.Font.Bold = 0
}

$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# Row 1, Col 2: "Sont les…" -> "Sont" / "…"
$c1 = $tbl.Cell(1, 2).Range
$c1.Find.Execute("Sont les…", $true, $false, $false, $false, $false, $true, 1, $false, "Sont…", 2)
Split-RunAt $d $c1 4

# Row 2, Col 2: "Télécom" / " Nancy" / "." -> single run "Télécom Nancy."
$c2 = $tbl.Cell(2, 2).Range
$c2.Find.Execute("Télécom Nancy.", $true, $false, $false, $false, $false, $true, 1, $false, "Télécom Nancy.", 2)

# Row 3, Col 2: "Les professeurs gérants les différents groupes" / "." ->
# single run "Les professeurs gérants les différents groupes."
$c3 = $tbl.Cell(3, 2).Range
$c3.Find.Execute("Les professeurs gérants les différents groupes.", $true, $false, $false, $false, $false, $true, 1, $false, "Les professeurs gérants les différents groupes.", 2)

# Row 4, Col 2: "Équipes projet." -> "Les élèves de l’é" / "quipe projet."
$c4 = $tbl.Cell(4, 2).Range
$c4.Find.Execute("Équipes projet.", $true, $false, $false, $false, $false, $true, 1, $false, "Les élèves de l’équipe projet.", 2)
Split-RunAt $d $c4 17

# Row 5, Col 2: "Communauté derrière les jardins partagés/circuit-court" ->
# "Les c" / "ommunauté" / "s" / " derrière les jardins partagés" / "/" /
# "circuit" / "s" / "-court" / "s."
$c5 = $tbl.Cell(5, 2).Range
$c5.Find.Execute("Communauté derrière les jardins partagés/circuit-court", $true, $false, $false, $false, $false, $true, 1, $false, "Les communautés derrière les jardins partagés/circuits-courts.", 2)
foreach ($off in @(5, 14, 15, 45, 46, 53, 54, 60)) {
    Split-RunAt $d $c5 $off
}
